# "finishing touches for deliverable 2"
# Fill in the previously-unknown ("TBD") Story Points estimates for each
# backlog item with their actual numeric values, and update the Total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = 8
$ws.Range("F3").Value  = 8
$ws.Range("F4").Value  = 6
$ws.Range("F5").Value  = 4
$ws.Range("F6").Value  = 4
$ws.Range("F7").Value  = 5
$ws.Range("F8").Value  = 2
$ws.Range("F9").Value  = 2
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 6
$ws.Range("F12").Value = 1

# Total row reflects the sum of the individual story points above.
$ws.Range("F18").Value = 48

# Leave the selection where the author last left it while wrapping up.
$ws.Range("F19").Select()
